$d = $word.ActiveDocument

$newParagraphs = @(
    @{ Text = "Folder: levelA";  Style = "Heading2" },
    @{ Text = "Folder: levelB1"; Style = "Heading3" },
    @{ Text = "Folder: levelB2"; Style = "Heading3" },
    @{ Text = "Folder: levelC";  Style = "Heading2" }
)

foreach ($item in $newParagraphs) {
    $endRange = $d.Content
    $endRange.Collapse(0)
    $endRange.InsertParagraphAfter()
    $endRange.Collapse(0)

    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $item.Text
    $newPara.Style = $item.Style
}

Write-Output ("Appended " + $newParagraphs.Count + " paragraphs; total paragraphs now " + $d.Paragraphs.Count)
